$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the formulas in C2:C4 (shared "1000/60" formula) with plain
# literal values of 1000 - the formulas are removed entirely.
$ws.Range("C2").Value = 1000
$ws.Range("C3").Value = 1000
$ws.Range("C4").Value = 1000

# Update the active selection on the sheet to a single cell, H2.
$ws.Range("H2").Select() | Out-Null
